$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add two new rows of data to Sheet1
$ws1.Range("A3").Value = 5
$ws1.Range("B3").Value = $false
$ws1.Range("C3").Value = """And Dad"""
$ws1.Range("D3").Value = 3

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = $true
$ws1.Range("C4").Value = """foobar"""
$ws1.Range("D4").Value = 6

# Update selection on Sheet1
$ws1.Range("E11").Select()

# Update selection on Sheet2 stays as A2, just unselected tab
$ws2.Range("A2").Select()

# Make Sheet1 the active sheet (so it is saved as the active tab)
$ws1.Activate()
$ws1.Range("E11").Select()
